$wb = $excel.ActiveWorkbook
$wsCat = $wb.Worksheets.Item(1)
$wsDatos = $wb.Worksheets.Item(2)

# 1. Remove the placeholder "aaaa" entry from the "tipos" list (Tabla2, datos!A4)
#    and shrink that table from A1:A4 to A1:A3.
$loTipos = $wsDatos.ListObjects.Item("Tabla2")
$wsDatos.Cells.Item(4, 1).ClearContents()
$loTipos.Resize($wsDatos.Range("A1:A3"))

# 2. Add the new "Pontevedra" category to the "categoría" list (Tabla3, datos!C1:C5 -> C1:C6)
$loCategoria = $wsDatos.ListObjects.Item("Tabla3")
$loCategoria.Resize($wsDatos.Range("C1:C6"))
$wsDatos.Cells.Item(6, 3).Value = "Pontevedra"

# 3. Add the new product row to the catalogo table (Tabla1, A1:I4 -> A1:I5)
$loProductos = $wsCat.ListObjects.Item("Tabla1")
$loProductos.Resize($wsCat.Range("A1:I5"))
$wsCat.Cells.Item(5, 1).Value = "Pontevedra Rafia"
$wsCat.Cells.Item(5, 2).Value = "bolsos"
$wsCat.Cells.Item(5, 3).Value = "Pontevedra"
$wsCat.Cells.Item(5, 4).Value = 18
$wsCat.Cells.Item(5, 4).NumberFormat = "#,##0\ ""€"";[Red]\-#,##0\ ""€"""
$wsCat.Cells.Item(5, 5).Value = "IMG_9458.HEIC"

# 4. Extend the data validation ranges on the catalogo sheet to cover the new row.
$wsCat.Range("B2:B5").Validation.Delete()
$wsCat.Range("B2:B5").Validation.Add(3, 1, 1, "=tipos")
$wsCat.Range("C2:C5").Validation.Delete()
$wsCat.Range("C2:C5").Validation.Add(3, 1, 1, "=categoria")

# 5. Update selection to match the saved workbook state.
$wsCat.Range("F9").Select() | Out-Null
